# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AC, AD, AE, styled like the other
# header cells (copy format from the existing header AB1 so the new
# cells share the same style index instead of minting a new one).
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-35: every row gets the same team record for this season.
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 104   # AC -> Wins
    $ws.Cells.Item($r, 30).Value = 58    # AD -> Losses
    $ws.Cells.Item($r, 31).Value = 0     # AE -> Ties
}
